# Apply updated cryptocurrency market data to the sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "27.834.20", "0.997").
# Pre-format the cells that will receive such values as Text so Excel
# stores them verbatim instead of coercing them into numbers/dates.
$ws.Range('D2:D5').NumberFormat = '@'
$ws.Range('D11:D13').NumberFormat = '@'
$ws.Range('D15:D19').NumberFormat = '@'
$ws.Range('D25:D26').NumberFormat = '@'
$ws.Range('D28:D29').NumberFormat = '@'
$ws.Range('D33:D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D39:D51').NumberFormat = '@'

$ws.Range('D2').Value = '27.834.20'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '1.628.97'
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '211.21'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('D11').Value = '0.0880'
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('D12').Value = '1.860.40'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').Value = '1.632.91'
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('E14').Value = '  -1.19%  '
$ws.Range('D15').Value = '0.555'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').Value = '64.81'
$ws.Range('E16').Value = '  -1.38%  '
$ws.Range('D17').Value = '27.852.80'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').Value = '228.27'
$ws.Range('E18').Value = '  -1.71%  '
$ws.Range('D19').Value = '7.61'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('E20').Value = '  -1.06%  '
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('E23').Value = '  -4.75%  '
$ws.Range('E24').Value = '  -0.55%  '
$ws.Range('D25').Value = '155.16'
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').Value = '6.91'
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('E27').Value = '  -0.28%  '
$ws.Range('D28').Value = '15.45'
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.414.79'
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '3.10'
$ws.Range('E34').Value = '  -0.83%  '
$ws.Range('E35').Value = '  +2.50%  '
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -3.76%  '
$ws.Range('E37').Value = '  -1.36%  '
$ws.Range('E38').Value = '  -1.13%  '
$ws.Range('D39').Value = '0.552'
$ws.Range('E39').Value = '  -1.00%  '
$ws.Range('D40').Value = '0.852'
$ws.Range('E40').Value = '  -1.79%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -1.89%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '65.66'
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '1.81'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').Value = '5.42'
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.769.56'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').Value = '2.12'
$ws.Range('E46').Value = '  -3.74%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '88.50'
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '0.101'
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0503'
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.61'
$ws.Range('E50').Value = '  +1.06%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '0.996'
$ws.Range('E51').Value = '  -0.21%  '
